$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 0, 103014),
    @(1, 0.2056551268874797, 104182),
    @(0.7871077816753987, 0.2056551268874797, 144902),
    @(0.7871077816753987, 0.2953213847735023, 146070),
    @(0.5458078442389603, 0.4757323837457365, 147238),
    @(0.5458078442389603, 0.4757323837457365, 186790),
    @(0.5458078442389603, 0.4757323837457365, 187958),
    @(0.5458078442389603, 0.5691517309668779, 189126),
    @(0.5458078442389603, 0.7140358456794925, 190294),
    @(0.5125494866107783, 0.7140358456794925, 228678),
    @(0.359736017481359, 0.7140358456794925, 229846),
    @(0.359736017481359, 0.7140358456794925, 231014),
    @(0.359736017481359, 0.7280207055893819, 232182),
    @(0.1620043134970038, 0.8252441279152454, 233350),
    @(0.1620043134970038, 0.8252441279152454, 270566),
    @(0.1620043134970038, 0.8252441279152454, 271734),
    @(0.1620043134970038, 0.8252441279152454, 272902),
    @(0.1620043134970038, 0.8252441279152454, 274070),
    @(0.1620043134970038, 0.8700772568582568, 275238),
    @(0.1020690647786086, 0.9682262050754993, 276406),
    @(0.1020690647786086, 0.9682262050754993, 312454),
    @(0.09408290090872004, 0.9682262050754993, 313622),
    @(0.05307277450263438, 0.9682262050754993, 314790),
    @(0.05307277450263438, 0.9682262050754993, 315958),
    @(0.03763316036348742, 0.9682262050754993, 317126),
    @(0.02246601727366483, 0.9682262050754993, 318294),
    @(0, 1, 319462)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
}
